$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "Male"
$ws.Range("C15").Value = 21
$ws.Range("D15").Value = "Student"
$ws.Range("E15").Value = 28
$ws.Range("F15").Value = "Excellent"
$ws.Range("G15").Value = "Work related"
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = "Never"
$ws.Range("L15").Value = "No"
$ws.Range("N15").Value = "Yes"
$ws.Range("O15").Value = "Samsung Notebook 10.1 Pro"
$ws.Range("P15").Value = "All the time"

$ws.Range("P15").Select()
